# Update the date in A1 and the price list in column D of Hoja1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date (A1) - serial 45436 = 2024-05-24
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Update prices in column D
$ws.Range("D29").Value = 2137
$ws.Range("D30").Value = 2322
$ws.Range("D31").Value = 1638
$ws.Range("D32").Value = 1726
$ws.Range("D33").Value = 1750
$ws.Range("D34").Value = 1861
$ws.Range("D35").Value = 1876
$ws.Range("D36").Value = 2052
